$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 134.83333
$ws.Range("I4").Value = 132
$ws.Range("K4").Value = 132
$ws.Range("M4").Value = -18

$ws.Range("H6").Value = 185.25
$ws.Range("I6").Value = 185.25
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 555.75
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -443.75
$ws.Range("N6").ClearContents()

$ws.Range("H12").Value = 622.63635
$ws.Range("I12").Value = 983.1667
$ws.Range("J12").Value = 190
$ws.Range("K12").Value = 983.1667
$ws.Range("L12").Value = 190
$ws.Range("M12").Value = -813.1667
$ws.Range("N12").Value = -530

$ws.Range("H19").Value = 653.0625
$ws.Range("I19").Value = 160
$ws.Range("J19").Value = 723.5
$ws.Range("K19").Value = 160
$ws.Range("L19").Value = 723.5
$ws.Range("M19").Value = 15
$ws.Range("N19").Value = -1073.5

$ws.Range("H53").Value = 139.38461
$ws.Range("I53").Value = 119.8
$ws.Range("K53").Value = 119.8
$ws.Range("M53").Value = 517.2

$ws.Range("H70").Value = 1002
$ws.Range("I70").Value = 1002
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3006
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2736
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 1002
$ws.Range("I73").Value = 1002
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3006
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2070
$ws.Range("N73").ClearContents()

$ws.Range("H115").Value = 789
$ws.Range("I115").Value = 861.25
$ws.Range("K115").Value = 2583.75
$ws.Range("M115").Value = -1016.75

$ws.Range("H129").Value = 4541.3335
$ws.Range("I129").Value = 5953.5
$ws.Range("J129").Value = 1717
$ws.Range("K129").Value = 17860.5
$ws.Range("L129").Value = 5151
$ws.Range("M129").Value = -12860.5
$ws.Range("N129").Value = -15151

$ws.Range("H135").Value = 1077
$ws.Range("I135").Value = 1077
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9693
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7158
$ws.Range("N135").ClearContents()

$ws.Range("H138").Value = 1989.3636
$ws.Range("I138").Value = 841
$ws.Range("K138").Value = 2523
$ws.Range("M138").Value = 2617

$ws.Range("H141").Value = 4268.125
$ws.Range("I141").Value = 4268.125
$ws.Range("K141").Value = 12804.375
$ws.Range("M141").Value = -7624.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H88").Value = 3112.7144
$ws.Range("J88").Value = 3299
$ws.Range("L88").Value = 3299
$ws.Range("N88").Value = -4111

$ws.Range("H91").Value = 3112.7144
$ws.Range("J91").Value = 3299
$ws.Range("L91").Value = 3299
$ws.Range("N91").Value = -6107

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1233
$ws.Range("I20").Value = 1233
$ws.Range("K20").Value = 1233
$ws.Range("M20").Value = -986

$ws.Range("H94").Value = 2428.4285
$ws.Range("I94").Value = 1928.4286
$ws.Range("K94").Value = 1928.4286
$ws.Range("M94").Value = -1477.4286

$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 600
$ws.Range("M107").Value = 1320

$ws.Range("H108").Value = 50000000
$ws.Range("J108").Value = 50000000
$ws.Range("L108").Value = 50000000
$ws.Range("N108").Value = -50007680

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 780
$ws.Range("I22").Value = 816.4
$ws.Range("J22").Value = 598
$ws.Range("K22").Value = 816.4
$ws.Range("L22").Value = 598
$ws.Range("M22").Value = -466.4
$ws.Range("N22").Value = -1298

$ws.Range("H86").Value = 8368.875
$ws.Range("I86").Value = 8571.714
$ws.Range("K86").Value = 8571.714
$ws.Range("M86").Value = -7448.714

$ws.Range("H88").Value = 27299.334
$ws.Range("J88").Value = 27299.334
$ws.Range("L88").Value = 27299.334
$ws.Range("N88").Value = -28111.334

$ws.Range("H89").Value = 8368.875
$ws.Range("I89").Value = 8571.714
$ws.Range("K89").Value = 42858.57
$ws.Range("M89").Value = -37242.57

$ws.Range("H91").Value = 27299.334
$ws.Range("J91").Value = 27299.334
$ws.Range("L91").Value = 27299.334
$ws.Range("N91").Value = -30107.334

$ws.Range("H132").Value = 1333
$ws.Range("I132").Value = 1333
$ws.Range("K132").Value = 3999
$ws.Range("M132").Value = -1469

$ws.Range("H134").Value = 1944.7142
$ws.Range("I134").Value = 1435.5
$ws.Range("K134").Value = 4306.5
$ws.Range("M134").Value = -1771.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 20.678572
$ws.Range("J2").Value = 32.25
$ws.Range("L2").Value = 193.5
$ws.Range("N2").Value = -419.5

$ws.Range("H38").Value = 760
$ws.Range("I38").Value = 108.333336
$ws.Range("K38").Value = 325.000008
$ws.Range("M38").Value = 21.99999200000002

$ws.Range("H80").Value = 1860
$ws.Range("J80").Value = 2100.6667
$ws.Range("L80").Value = 6302.000100000001
$ws.Range("N80").Value = -8174.000100000001

$ws.Range("H83").Value = 1860
$ws.Range("J83").Value = 2100.6667
$ws.Range("L83").Value = 18906.0003
$ws.Range("N83").Value = -28266.0003

$ws.Range("H134").Value = 992
$ws.Range("I134").Value = 992
$ws.Range("K134").Value = 2976
$ws.Range("M134").Value = 2094

$ws.Range("H139").Value = 2203.3333
$ws.Range("I139").Value = 2203.3333
$ws.Range("K139").Value = 6609.999899999999
$ws.Range("M139").Value = -1469.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6250000
$ws.Range("I3").Value = 5000000
$ws.Range("K3").Value = 5000000
$ws.Range("M3").Value = -4999884

$ws.Range("H11").Value = 4112500
$ws.Range("I11").Value = 1375000
$ws.Range("J11").Value = 6850000
$ws.Range("K11").Value = 1375000
$ws.Range("L11").Value = 6850000
$ws.Range("M11").Value = -1374861
$ws.Range("N11").Value = -6850278

$ws.Range("H33").Value = 5012500
$ws.Range("J33").Value = 5012500
$ws.Range("L33").Value = 5012500
$ws.Range("N33").Value = -5013004

$ws.Range("H82").Value = 200328
$ws.Range("J82").Value = 200328
$ws.Range("L82").Value = 200328
$ws.Range("N82").Value = -201094

$ws.Range("H85").Value = 200328
$ws.Range("J85").Value = 200328
$ws.Range("L85").Value = 200328
$ws.Range("N85").Value = -202980

$ws.Range("H102").Value = 25958.375
$ws.Range("I102").Value = 25958.375
$ws.Range("K102").Value = 25958.375
$ws.Range("M102").Value = -24336.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5818.2
$ws.Range("I7").Value = 5999
$ws.Range("J7").Value = 5396.3335
$ws.Range("K7").Value = 5999
$ws.Range("L7").Value = 5396.3335
$ws.Range("M7").Value = -5887
$ws.Range("N7").Value = -5620.3335

$ws.Range("H22").Value = 3161.875
$ws.Range("I22").Value = 2799
$ws.Range("J22").Value = 3213.7144
$ws.Range("K22").Value = 2799
$ws.Range("L22").Value = 3213.7144
$ws.Range("M22").Value = -2504
$ws.Range("N22").Value = -3803.7144

$ws.Range("H27").Value = 3161.875
$ws.Range("I27").Value = 2799
$ws.Range("J27").Value = 3213.7144
$ws.Range("K27").Value = 2799
$ws.Range("L27").Value = 3213.7144
$ws.Range("M27").Value = -2692
$ws.Range("N27").Value = -3427.7144

$ws.Range("H55").Value = 464.17648
$ws.Range("I55").Value = 334.75
$ws.Range("J55").Value = 504
$ws.Range("K55").Value = 334.75
$ws.Range("L55").Value = 504
$ws.Range("M55").Value = -161.75
$ws.Range("N55").Value = -850

$ws.Range("H122").Value = 2319.3333
$ws.Range("I122").Value = 1984.25
$ws.Range("K122").Value = 5952.75
$ws.Range("M122").Value = -3502.75

$ws.Range("H126").Value = 5818.2
$ws.Range("I126").Value = 5999
$ws.Range("J126").Value = 5396.3335
$ws.Range("K126").Value = 17997
$ws.Range("L126").Value = 16189.0005
$ws.Range("M126").Value = -15527
$ws.Range("N126").Value = -21129.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 4500
$ws.Range("M132").Value = -1970
